$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1) ---
$ws.Cells.Item(1,1).Value = "Date"
$ws.Cells.Item(1,2).Value = "ScoreFinal"
$ws.Cells.Item(1,3).Value = "totalSentiment"
$ws.Cells.Item(1,4).Value = "posWordPercentage"
$ws.Cells.Item(1,5).Value = "negWordPercentage"
$ws.Cells.Item(1,6).Value = "posPhrasePercentage"
$ws.Cells.Item(1,7).Value = "negPhrasePercentage"
$ws.Cells.Item(1,8).Value = "ElapsedMs"
$ws.Cells.Item(1,9).Value = "wordCount"
$ws.Cells.Item(1,10).Value = "sentenceCount"
$ws.Cells.Item(1,11).Value = "posWordCount"
$ws.Cells.Item(1,12).Value = "negWordCount"
$ws.Cells.Item(1,13).Value = "positivePhraseCount"
$ws.Cells.Item(1,14).Value = "negativePhraseCount"
$ws.Cells.Item(1,15).Value = "Method"
$ws.Cells.Item(1,16).Value = "RSI"
$ws.Cells.Item(1,17).Value = "PEG"

# --- Row 2 ---
$ws.Cells.Item(2,1).Value = 42627.874236111114
$ws.Cells.Item(2,2).Value = 8
$ws.Cells.Item(2,3).Value = 40
$ws.Cells.Item(2,4).Value = 65
$ws.Cells.Item(2,5).Value = 33
$ws.Cells.Item(2,6).Value = 99
$ws.Cells.Item(2,7).Value = 0
$ws.Cells.Item(2,8).Value = 9555
$ws.Cells.Item(2,9).Value = 4917
$ws.Cells.Item(2,10).Value = 743
$ws.Cells.Item(2,11).Value = 122
$ws.Cells.Item(2,12).Value = 62
$ws.Cells.Item(2,13).Value = 12
$ws.Cells.Item(2,14).Value = 0
$ws.Cells.Item(2,15).Value = "Noun"
$ws.Cells.Item(2,16).Value = 0
$ws.Cells.Item(2,17).Value = 1

# --- Row 3 ---
$ws.Cells.Item(3,1).Value = 42627.87724537037
$ws.Cells.Item(3,2).Value = 4
$ws.Cells.Item(3,3).Value = 22
$ws.Cells.Item(3,4).Value = 64
$ws.Cells.Item(3,5).Value = 35
$ws.Cells.Item(3,6).Value = 100
$ws.Cells.Item(3,7).Value = 0
$ws.Cells.Item(3,8).Value = 7388
$ws.Cells.Item(3,9).Value = 3830
$ws.Cells.Item(3,10).Value = 545
$ws.Cells.Item(3,11).Value = 81
$ws.Cells.Item(3,12).Value = 44
$ws.Cells.Item(3,13).Value = 1
$ws.Cells.Item(3,14).Value = 0
$ws.Cells.Item(3,15).Value = "Noun"
$ws.Cells.Item(3,16).Value = 0
$ws.Cells.Item(3,17).Value = 1
